# Fruta / hortaliza, semanal
# Insert a new weekly record at row 925, pushing the existing rows
# (925-961) down to (926-962).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 925; this shifts rows
# 925..961 down to 926..962 and extends the sheet dimension to R962.
$ws.Rows("925:925").Insert()

# Populate the newly inserted row 925 with the new data point.
$ws.Cells.Item(925, 1).Value = 3
$ws.Cells.Item(925, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(925, 3).Value = "Coquimbo"
$ws.Cells.Item(925, 4).Value = 45147
$ws.Cells.Item(925, 5).Value = 5
$ws.Cells.Item(925, 6).Value = 100112045
$ws.Cells.Item(925, 7).Value = "Zapallo"
$ws.Cells.Item(925, 8).Value = "Camote"
$ws.Cells.Item(925, 9).Value = "1a (guarda)"
$ws.Cells.Item(925, 10).Value = 120
$ws.Cells.Item(925, 11).Value = 500
$ws.Cells.Item(925, 12).Value = 500
$ws.Cells.Item(925, 13).Value = 500
$ws.Cells.Item(925, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(925, 15).Value = "Provincia de Talca"
$ws.Cells.Item(925, 16).Value = 500
$ws.Cells.Item(925, 17).Value = 1
$ws.Cells.Item(925, 18).Value = "Hortaliza"
